$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Messages")

# --- New "Slot number" column (J) alongside the existing message rows 2-6 ---
# Fill used: light gray (theme "Background 1, Darker 15%" -> RGB D9D9D9)
$grayFill = 14277081   # RGB(217,217,217) == theme 0 (Background1) tinted -15%

$ws.Range("J2").Value = "Slot number"
$ws.Range("J2").Interior.Color = $grayFill

$ws.Range("J3").Value = "Slot number"
$ws.Range("J3").Interior.Color = $grayFill

$ws.Range("J4").Value = "Slot number"
$ws.Range("J4").Interior.Color = $grayFill

# --- Row 5: insert a new "Slot time" cell in C5, push the old "Desired state"
#     value (with its original formatting) one column over, into D5 ---
$ws.Range("C5").Copy()
$ws.Range("D5").PasteSpecial(-4122)              # xlPasteFormats
$ws.Range("D5").Value = $ws.Range("C5").Value2   # "Desired state"

$ws.Range("C5").Value = "Slot time"
$ws.Range("C5").Interior.ThemeColor = 6          # theme 5 (Accent2), no tint

$ws.Range("J5").Value = "Slot number"
$ws.Range("J5").Interior.Color = $grayFill

$ws.Range("J6").Value = "Slot number"
$ws.Range("J6").Interior.Color = $grayFill

# --- Selection moved to C9 on the Messages sheet ---
$ws.Range("C9").Select()
